$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare numeric literal (e.g. "223.50", "1.00",
# "0.0000292") would otherwise be auto-coerced to a Number by Excel and lose
# the exact original text (trailing zeros, dot-grouping, sci-notation-prone
# tiny decimals). Force those specific cells to Text before writing so the
# literal string is preserved, matching the source data (inline strings).

$ws.Range("D2").Value = '95.702.77'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '3.622.73'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '2.71'
$ws.Range("E4").Value = '  +38.04%  '
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '223.50'
$ws.Range("E6").Value = '  -5.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '639.45'
$ws.Range("E7").Value = '  -2.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.421'
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("E9").Value = '  +12.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("D11").Value = '3.619.88'
$ws.Range("E11").Value = '  -2.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.24'
$ws.Range("E12").Value = '  +7.93%  '
$ws.Range("E13").Value = '  +2.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000292'
$ws.Range("E14").Value = '  -6.83%  '
$ws.Range("E15").Value = '  -5.71%  '
$ws.Range("D16").Value = '4.297.10'
$ws.Range("E16").Value = '  -2.58%  '
$ws.Range("D17").Value = '95.459.20'
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '22.90'
$ws.Range("E18").Value = '  +22.38%  '
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.90'
$ws.Range("E20").Value = '  +7.03%  '
$ws.Range("D21").Value = '3.617.67'
$ws.Range("E21").Value = '  -2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.296'
$ws.Range("E22").Value = '  +52.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.542'
$ws.Range("E23").Value = '  +5.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '515.81'
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("E25").Value = '  -6.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '125.91'
$ws.Range("E26").Value = '  +18.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000201'
$ws.Range("E27").Value = '  -10.56%  '
$ws.Range("E28").Value = '  -1.25%  '
$ws.Range("D29").Value = '3.805.87'
$ws.Range("E29").Value = '  -2.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.74'
$ws.Range("E30").Value = '  -5.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.05'
$ws.Range("E31").Value = '  +3.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.07'
$ws.Range("E32").Value = '  +1.20%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("E34").Value = '  +4.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.182'
$ws.Range("E35").Value = '  -5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.70'
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.539'
$ws.Range("E39").Value = '  +6.87%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.20'
$ws.Range("E41").Value = '  +6.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '583.55'
$ws.Range("E42").Value = '  -9.23%  '
$ws.Range("E43").Value = '  -4.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0528'
$ws.Range("E44").Value = '  +15.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.90'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.962'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("E47").Value = '  -5.72%  '
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.09'
$ws.Range("E49").Value = '  +4.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '233.39'
$ws.Range("E50").Value = '  +13.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.51'
$ws.Range("E51").Value = '  -0.42%  '
